$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = "Match (f)<--(g:genomic_info)`nWHERE g.library_selection in ['Hybrid Selection']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p, s, apoc.coll.sort(collect(distinct(samp.sample_id))) as samples`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samples, ','), '') as ``Samples```nORDER BY ``Participant ID``LIMIT 100"

$ws.Range("B2").Value = $newQuery

$ws.Rows.Item(2).RowHeight = 382.5

$ws.Range("B2").Select()
